# Update met oude projecten
# Adds a new column AO "GTA.GTA_PEILBUISGEGEVENS.PEILBUISIDENT" with
# PEILBUISIDENT values for each data row (rows 2-28), matching the
# header style already used by the other header cells in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = [ordered]@{
    1  = 'GTA.GTA_PEILBUISGEGEVENS.PEILBUISIDENT'
    2  = 'PB_1_001_04P001442_F-374'
    3  = 'PB_1_001_04P001442_F-924'
    4  = 'PB_1_001_04P001442_F-2224'
    5  = 'PB_1_002_04P001442_F-375'
    6  = 'PB_1_002_04P001442_F-975'
    7  = 'PB_1_002_04P001442_F-2525'
    8  = 'PB_1_003_04P001442_F-337'
    9  = 'PB_1_003_04P001442_F-937'
    10 = 'PB_1_003_04P001442_F-2987'
    11 = 'PB_1_004_04P001442-01_F-204'
    12 = 'PB_1_004_04P001442-01_F-1224'
    13 = 'PB_1_004_04P001442-01_F-2274'
    14 = 'PB_1_005_04P001442-01_F-264'
    15 = 'PB_1_005_04P001442-01_F-964'
    16 = 'PB_1_005_04P001442-01_F-2289'
    17 = 'PB_2_001_04P001442_F-209'
    18 = 'PB_2_001_04P001442_F-889'
    19 = 'PB_2_001_04P001442_F-2889'
    20 = 'PB_2_002_04P001442_F-297'
    21 = 'PB_2_002_04P001442_F-897'
    22 = 'PB_2_002_04P001442_F-2497'
    23 = 'PB_2_003_04P001442-01_F-206'
    24 = 'PB_2_003_04P001442-01_F-906'
    25 = 'PB_2_003_04P001442-01_F-2356'
    26 = 'PB_6_001_04P001442-01_F-245'
    27 = 'PB_6_001_04P001442-01_F-875'
    28 = 'PB_6_001_04P001442-01_F-2325'
}

# Column AO is column index 41
$col = 41

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, $col).Value = $values[$row]
}

# Header cell AO1 uses the same formatting as the other header cells (e.g. AN1)
$ws.Range("AN1").Copy()
$ws.Range("AO1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
